$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Examples & Info")

# Row 2 (header descriptions) - fill in real text where placeholder
# "another description" text existed, and fix a typo.
$ws.Range("Q2").Value  = "How was the mouse/animal killed?"
$ws.Range("R2").Value  = "Any treatment / pertubation applied to the individual/ cell-line."
$ws.Range("AF2").Value = "Nucleic acid concentration"
$ws.Range("AG2").Value = "Average fragment length (basepairs)"
$ws.Range("AH2").Value = "Library molarity"
$ws.Range("AJ2").Value = "Well position of the barcode"
$ws.Range("AL2").Value = "Well position of the barcode"
$ws.Range("AV2").Value = "Location of the library/experiement description"
$ws.Range("AW2").Value = "Plate label used"

# Row 6 (regex examples row) - add date-format regex examples
$ws.Range("O6").Value  = "[0-9]{4}-[0-9]{2}-[0-9]{2}"
$ws.Range("P6").Value  = "[0-9]{4}-[0-9]{2}-[0-9]{2}"
$ws.Range("AU6").Value = "[0-9]{4}-[0-9]{2}-[0-9]{2}"
